$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200; this shifts the existing rows
# 200-290 down to 201-291 (matching dimension growing to A1:R291).
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new record's data.
$ws.Range("A200").Value = 3
$ws.Range("B200").Value = "Femacal de La Calera"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44489
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = 100112021
$ws.Range("G200").Value = "Ají"
$ws.Range("H200").Value = "Americana (o)"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 45
$ws.Range("K200").Value = 40000
$ws.Range("L200").Value = 41000
$ws.Range("M200").Value = 40444
$ws.Range("N200").Value = "$/caja 15 kilos"
$ws.Range("O200").Value = "Región de Arica y Parinacota"
$ws.Range("P200").Value = 2696
$ws.Range("Q200").Value = 15
$ws.Range("R200").Value = "Hortaliza"
